$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.335.61"
$ws.Range("E2").Value = "  +2.27%  "

# Row 3
$ws.Range("D3").Value = "2.692.33"
$ws.Range("E3").Value = "  +2.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.39%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.29"
$ws.Range("E5").Value = "  +2.84%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.22"
$ws.Range("E6").Value = "  +3.89%  "

# Row 7
$ws.Range("E7").Value = "  -0.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +0.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.127"
$ws.Range("E9").Value = "  +10.36%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.01"
$ws.Range("E10").Value = "  +3.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.408"
$ws.Range("E11").Value = "  +2.81%  "

# Row 12
$ws.Range("E12").Value = "  +1.75%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  +25.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.42"
$ws.Range("E14").Value = "  +5.31%  "

# Row 15
$ws.Range("D15").Value = "3.173.46"
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$ws.Range("D16").Value = "66.137.63"
$ws.Range("E16").Value = "  +2.13%  "

# Row 17
$ws.Range("D17").Value = "2.679.38"
$ws.Range("E17").Value = "  +2.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.79"
$ws.Range("E18").Value = "  +3.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.93"
$ws.Range("E19").Value = "  +2.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.26"
$ws.Range("E20").Value = "  +3.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.51"
$ws.Range("E21").Value = "  +4.66%  "

# Row 22
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.10"
$ws.Range("E23").Value = "  +3.51%  "

# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000109"
$ws.Range("E24").Value = "  +19.06%  "

# Row 25
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.70"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.67"
$ws.Range("E26").Value = "  +3.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.174"
$ws.Range("E27").Value = "  +5.46%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.66"
$ws.Range("E28").Value = "  +0.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("E29").Value = "  -0.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +7.76%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "543.97"
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.59"
$ws.Range("E34").Value = "  +5.34%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.61"
$ws.Range("E35").Value = "  -1.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.437"
$ws.Range("E36").Value = "  +3.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.86"
$ws.Range("E37").Value = "  +3.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.03"
$ws.Range("E38").Value = "  +0.60%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.81"
$ws.Range("E39").Value = "  -1.11%  "

# Row 40
$ws.Range("E40").Value = "  -0.31%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.70"
$ws.Range("E41").Value = "  +2.26%  "

# Row 42
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.43"
$ws.Range("E43").Value = "  +1.33%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.22"
$ws.Range("E44").Value = "  +3.13%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("E45").Value = "  +7.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0621"
$ws.Range("E46").Value = "  +2.94%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.34"
$ws.Range("E47").Value = "  -0.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.665"
$ws.Range("E48").Value = "  +3.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0266"
$ws.Range("E49").Value = "  +5.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.07"
$ws.Range("E50").Value = "  +3.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0986"
$ws.Range("E51").Value = "  +0.28%  "
